$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the three previously-blank rows (9, 10, 11) with new note-related
# error entries (columns: A=Code, B=Reason, C=Error Message, D=Where, E=HTTP)

$ws.Range("A9").Value = "N0403"
$ws.Range("B9").Value = "User is updating a note which does not exist."
$ws.Range("C9").Value = "This note does not exist."
$ws.Range("D9").Value = "backend.notes.views.updateNote"
$ws.Range("E9").Value = 404

$ws.Range("A10").Value = "N0404"
$ws.Range("B10").Value = "The user does not have or own that note."
$ws.Range("C10").Value = "This note does not exist."
$ws.Range("D10").Value = "backend.notes.views.readNote"
$ws.Range("E10").Value = 404

$ws.Range("A11").Value = "N0407"
$ws.Range("B11").Value = "User is trying to delete a note which does not exist or they do not own."
$ws.Range("C11").Value = "This note does not exist."
$ws.Range("D11").Value = "backend.notes.views.deleteNote"
$ws.Range("E11").Value = 404

# Append two more blank placeholder rows at the bottom (59, 60), matching the
# formatting of the other blank rows in column A (style carried from A8)
$ws.Range("A8").Copy()
$ws.Range("A59:A60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection to reflect where the author left off editing
$ws.Range("C24").Select()
